$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K" - strikeouts) values per row, replacing old Strike# figures
$gValues = @{
    2 = 2
    3 = 3
    4 = 0
    5 = 2
    6 = 3
    7 = 1
    8 = 1
    9 = 1
    10 = 2
    11 = 1
    12 = 2
    13 = 2
    14 = 1
    15 = 0
    16 = 1
    17 = 0
    19 = 0
    20 = 2
    21 = 0
    22 = 2
    23 = 1
    24 = 1
    25 = 0
    26 = 2
    27 = 1
    28 = 1
    29 = 1
    30 = 0
    31 = 3
    32 = 1
    33 = 0
    34 = 0
    35 = 1
    36 = 1
    37 = 3
    38 = 1
    39 = 3
    40 = 2
    41 = 2
    42 = 2
    43 = 3
    44 = 0
    45 = 1
    46 = 0
    47 = 1
    48 = 2
    49 = 1
    50 = 2
    51 = 2
    52 = 1
    53 = 2
    54 = 2
    55 = 2
    56 = 1
    57 = 0
    58 = 1
    59 = 2
    60 = 1
    61 = 2
    62 = 0
    63 = 2
    64 = 2
    65 = 2
    66 = 0
    67 = 2
    68 = 3
    69 = 2
    70 = 1
    71 = 1
    72 = 1
    73 = 1
    74 = 1
    75 = 0
    76 = 1
    77 = 1
    78 = 0
    79 = 2
    82 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

